# Auto-generated-derived edit script for Hades_Profits workbook
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H-N) for specific rows
# across sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR based on the authoritative diff.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 74 (Leve Item ID 5507)
$ws.Range("H74").Value = 3912.64
$ws.Range("I74").Value = 3863.25
$ws.Range("J74").Value = 3935.8823
$ws.Range("K74").Value = 3863.25
$ws.Range("L74").Value = 3935.8823
$ws.Range("M74").Value = -2927.25
$ws.Range("N74").Value = -5807.8823
# Row 77 (Leve Item ID 5507)
$ws.Range("H77").Value = 3912.64
$ws.Range("I77").Value = 3863.25
$ws.Range("J77").Value = 3935.8823
$ws.Range("K77").Value = 19316.25
$ws.Range("L77").Value = 19679.4115
$ws.Range("M77").Value = -14636.25
$ws.Range("N77").Value = -29039.4115
# Row 98 (Leve Item ID 36237)
$ws.Range("H98").Value = 670.2593000000001
$ws.Range("I98").Value = 487.88
$ws.Range("J98").Value = 2950
$ws.Range("K98").Value = 487.88
$ws.Range("L98").Value = 2950
$ws.Range("M98").Value = 1010.12
$ws.Range("N98").Value = -5946
# Row 106 (Leve Item ID 19903)
$ws.Range("H106").Value = 2369
$ws.Range("I106").Value = 2377.6155
$ws.Range("J106").Value = 2331.6667
$ws.Range("K106").Value = 2377.6155
$ws.Range("L106").Value = 2331.6667
$ws.Range("M106").Value = -1746.6155
$ws.Range("N106").Value = -3593.6667
# Row 107 (Leve Item ID 27766)
$ws.Range("H107").Value = 519.86957
$ws.Range("I107").Value = 355.6154
$ws.Range("J107").Value = 733.4
$ws.Range("K107").Value = 355.6154
$ws.Range("L107").Value = 733.4
$ws.Range("M107").Value = 1564.3846
$ws.Range("N107").Value = -4573.4
# Row 113 (Leve Item ID 27775)
$ws.Range("H113").Value = 4710
$ws.Range("I113").Value = 4717.5
$ws.Range("J113").Value = 4701
$ws.Range("K113").Value = 4717.5
$ws.Range("L113").Value = 4701
$ws.Range("M113").Value = -1463.5
$ws.Range("N113").Value = -11209
# Row 114 (Leve Item ID 25959)
$ws.Range("H114").Value = 29888.4
$ws.Range("J114").Value = 29888.4
$ws.Range("L114").Value = 29888.4
$ws.Range("N114").Value = -38566.4
# Row 115 (Leve Item ID 27957)
$ws.Range("H115").Value = 1800
$ws.Range("I115").Value = 733.3333
$ws.Range("J115").Value = 5000
$ws.Range("K115").Value = 2199.9999
$ws.Range("L115").Value = 15000
$ws.Range("M115").Value = -632.9998999999998
$ws.Range("N115").Value = -18134
# Row 118 (Leve Item ID 27958)
$ws.Range("H118").Value = 2804.5
$ws.Range("I118").Value = 0
$ws.Range("J118").Value = 2804.5
$ws.Range("K118").Value = 0
$ws.Range("M118").Value = 8413.5
$ws.Range("N118").Value = -11727.5
# Row 120 (Leve Item ID 26279)
$ws.Range("H120").Value = 29480.5
$ws.Range("J120").Value = 29480.5
$ws.Range("L120").Value = 29480.5
$ws.Range("N120").Value = -39156.5
# Row 122 (Leve Item ID 36237)
$ws.Range("H122").Value = 670.2593000000001
$ws.Range("I122").Value = 487.88
$ws.Range("J122").Value = 2950
$ws.Range("K122").Value = 1463.64
$ws.Range("L122").Value = 8850
$ws.Range("M122").Value = 986.3600000000001
$ws.Range("N122").Value = -13750
# Row 123 (Leve Item ID 34090)
$ws.Range("H123").Value = 87636.37
$ws.Range("J123").Value = 87636.37
$ws.Range("L123").Value = 87636.37
$ws.Range("N123").Value = -97436.37
# Row 124 (Leve Item ID 34241)
$ws.Range("H124").Value = 35000
$ws.Range("J124").Value = 35000
$ws.Range("L124").Value = 35000
$ws.Range("N124").Value = -44820

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 2 (Leve Item ID 27713)
$ws.Range("H2").Value = 3960.2144
$ws.Range("I2").Value = 3804.889
$ws.Range("J2").Value = 4239.8
$ws.Range("K2").Value = 3804.889
$ws.Range("L2").Value = 4239.8
$ws.Range("M2").Value = -3691.889
$ws.Range("N2").Value = -4465.8
# Row 74 (Leve Item ID 44000)
$ws.Range("H74").Value = 11180829
$ws.Range("I74").Value = 15198693
$ws.Range("K74").Value = 15198693
$ws.Range("M74").Value = -15197819
# Row 77 (Leve Item ID 44000)
$ws.Range("H77").Value = 11180829
$ws.Range("I77").Value = 15198693
$ws.Range("K77").Value = 75993465
$ws.Range("M77").Value = -75989097
# Row 96 (Leve Item ID 18207)
$ws.Range("H96").Value = 44806
$ws.Range("J96").Value = 44806
$ws.Range("L96").Value = 44806
$ws.Range("N96").Value = -50298
# Row 116 (Leve Item ID 27713)
$ws.Range("H116").Value = 3960.2144
$ws.Range("I116").Value = 3804.889
$ws.Range("J116").Value = 4239.8
$ws.Range("K116").Value = 3804.889
$ws.Range("L116").Value = 4239.8
$ws.Range("M116").Value = -1510.889
$ws.Range("N116").Value = -8827.799999999999
# Row 117 (Leve Item ID 26125)
$ws.Range("H117").Value = 39800
$ws.Range("J117").Value = 39800
$ws.Range("L117").Value = 39800
$ws.Range("N117").Value = -48978
# Row 122 (Leve Item ID 36168)
$ws.Range("H122").Value = 15875989
$ws.Range("I122").Value = 4012
$ws.Range("J122").Value = 18521318
$ws.Range("K122").Value = 12036
$ws.Range("L122").Value = 55563954
$ws.Range("M122").Value = -9586
$ws.Range("N122").Value = -55568854
# Row 138 (Leve Item ID 42350)
$ws.Range("H138").Value = 86399.75
$ws.Range("J138").Value = 86399.75
$ws.Range("L138").Value = 86399.75
$ws.Range("N138").Value = -96679.75

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 3 (Leve Item ID 27713)
$ws.Range("H3").Value = 3960.2144
$ws.Range("I3").Value = 3804.889
$ws.Range("J3").Value = 4239.8
$ws.Range("K3").Value = 3804.889
$ws.Range("L3").Value = 4239.8
$ws.Range("M3").Value = -3690.889
$ws.Range("N3").Value = -4467.8

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 117 (Leve Item ID 27135)
$ws.Range("H117").Value = 39800
$ws.Range("J117").Value = 39800
$ws.Range("L117").Value = 39800
$ws.Range("N117").Value = -48978
# Row 141 (Leve Item ID 43345)
$ws.Range("H141").Value = 45386.55
$ws.Range("I141").Value = 22574
$ws.Range("J141").Value = 51089.688
$ws.Range("K141").Value = 22574
$ws.Range("L141").Value = 51089.688
$ws.Range("M141").Value = -17394
$ws.Range("N141").Value = -61449.688

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 7 (Leve Item ID 4728)
$ws.Range("H7").Value = 330.33334
$ws.Range("I7").Value = 197.38461
$ws.Range("J7").Value = 676
$ws.Range("K7").Value = 592.15383
$ws.Range("L7").Value = 2028
$ws.Range("M7").Value = -480.15383
$ws.Range("N7").Value = -2252
# Row 92 (Leve Item ID 19841)
$ws.Range("H92").Value = 1021.1724
$ws.Range("I92").Value = 1052.9565
$ws.Range("K92").Value = 3158.8695
$ws.Range("M92").Value = -1910.8695

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 113 (Leve Item ID 27710)
$ws.Range("H113").Value = 2316.3914
$ws.Range("I113").Value = 1589.5
$ws.Range("K113").Value = 1589.5
$ws.Range("M113").Value = 580.5
# Row 138 (Leve Item ID 42325)
$ws.Range("H138").Value = 56666
$ws.Range("J138").Value = 56666
$ws.Range("L138").Value = 56666
$ws.Range("N138").Value = -66946

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 46 (Leve Item ID 5282)
$ws.Range("H46").Value = 2165673.5
$ws.Range("I46").Value = 4329605.5
$ws.Range("J46").Value = 1741.1428
$ws.Range("K46").Value = 4329605.5
$ws.Range("L46").Value = 1741.1428
$ws.Range("M46").Value = -4329417.5
$ws.Range("N46").Value = -2117.1428
# Row 55 (Leve Item ID 5284)
$ws.Range("H55").Value = 241.3077
$ws.Range("I55").Value = 225.22728
$ws.Range("J55").Value = 329.75
$ws.Range("K55").Value = 225.22728
$ws.Range("L55").Value = 329.75
$ws.Range("M55").Value = -52.22728000000001
$ws.Range("N55").Value = -675.75
# Row 61 (Leve Item ID 27740)
$ws.Range("H61").Value = 1902.258
$ws.Range("I61").Value = 1898.0952
$ws.Range("J61").Value = 1911
$ws.Range("K61").Value = 1898.0952
$ws.Range("L61").Value = 1911
$ws.Range("M61").Value = -1696.0952
$ws.Range("N61").Value = -2315
# Row 113 (Leve Item ID 27740)
$ws.Range("H113").Value = 1902.258
$ws.Range("I113").Value = 1898.0952
$ws.Range("J113").Value = 1911
$ws.Range("K113").Value = 1898.0952
$ws.Range("L113").Value = 1911
$ws.Range("M113").Value = 271.9048
$ws.Range("N113").Value = -6251
# Row 132 (Leve Item ID 44058)
$ws.Range("H132").Value = 29173.816
$ws.Range("I132").Value = 2868.5217
$ws.Range("J132").Value = 69508.60000000001
$ws.Range("K132").Value = 8605.5651
$ws.Range("L132").Value = 208525.8
$ws.Range("M132").Value = -6075.5651
$ws.Range("N132").Value = -213585.8

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 81 (Leve Item ID 12596)
$ws.Range("H81").Value = 2066.6667
# Row 84 (Leve Item ID 12596)
$ws.Range("H84").Value = 2066.6667
# Row 113 (Leve Item ID 27752)
$ws.Range("H113").Value = 779.25
$ws.Range("I113").Value = 997.9231
$ws.Range("J113").Value = 589.73334
$ws.Range("K113").Value = 2993.7693
$ws.Range("L113").Value = 1769.20002
$ws.Range("M113").Value = -823.7692999999999
$ws.Range("N113").Value = -6109.20002
